$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 362.0435
$ws.Range("I19").Value = 351.73334
$ws.Range("J19").Value = 381.375
$ws.Range("K19").Value = 351.73334
$ws.Range("L19").Value = 381.375
$ws.Range("M19").Value = -176.73334
$ws.Range("N19").Value = -731.375

$ws.Range("H112").Value = 11628864
$ws.Range("I112").Value = 782.5
$ws.Range("J112").Value = 14286712
$ws.Range("K112").Value = 2347.5
$ws.Range("L112").Value = 42860136
$ws.Range("M112").Value = -1239.5
$ws.Range("N112").Value = -42862352

$ws.Range("H118").Value = 797
$ws.Range("I118").Value = 356
$ws.Range("J118").Value = 1127.75
$ws.Range("K118").Value = 1068
$ws.Range("L118").Value = 3383.25
$ws.Range("M118").Value = 589
$ws.Range("N118").Value = -6697.25

$ws.Range("H132").Value = 1755
$ws.Range("I132").Value = 1634.3226
$ws.Range("J132").Value = 3002
$ws.Range("K132").Value = 4902.9678
$ws.Range("L132").Value = 9006
$ws.Range("M132").Value = -2372.9678
$ws.Range("N132").Value = -14066

$ws.Range("H138").Value = 14546.158
$ws.Range("I138").Value = 1015.3684
$ws.Range("J138").Value = 17928.855
$ws.Range("K138").Value = 3046.1052
$ws.Range("L138").Value = 53786.565
$ws.Range("M138").Value = 2093.8948
$ws.Range("N138").Value = -64066.565

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4409.7417
$ws.Range("I32").Value = 3670.2927
$ws.Range("J32").Value = 5853.4287
$ws.Range("K32").Value = 3670.2927
$ws.Range("L32").Value = 5853.4287
$ws.Range("M32").Value = -3383.2927
$ws.Range("N32").Value = -6427.4287

$ws.Range("H45").Value = 10000
$ws.Range("I45").Value = 10000
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 10000
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -9623
$ws.Range("N45").ClearContents()

$ws.Range("H122").Value = 1710976.5
$ws.Range("I122").Value = 2138262.2
$ws.Range("J122").Value = 1833.3334
$ws.Range("K122").Value = 6414786.600000001
$ws.Range("L122").Value = 5500.0002
$ws.Range("M122").Value = -6412336.600000001
$ws.Range("N122").Value = -10400.0002

$ws.Range("H132").Value = 4985.8696
$ws.Range("I132").Value = 1914.2084
$ws.Range("J132").Value = 8336.772000000001
$ws.Range("K132").Value = 5742.6252
$ws.Range("L132").Value = 25010.316
$ws.Range("M132").Value = -3212.6252
$ws.Range("N132").Value = -30070.316

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1972.0834
$ws.Range("I86").Value = 1972.0834
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 1972.0834
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -849.0834
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 1972.0834
$ws.Range("I89").Value = 1972.0834
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 9860.416999999999
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -4244.416999999999
$ws.Range("N89").ClearContents()

$ws.Range("H134").Value = 3554.3333
$ws.Range("I134").Value = 3866.186
$ws.Range("J134").Value = 2765.5293
$ws.Range("K134").Value = 11598.558
$ws.Range("L134").Value = 8296.5879
$ws.Range("M134").Value = -9063.558000000001
$ws.Range("N134").Value = -13366.5879

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9930.286
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 9930.286
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 9930.286
$ws.Range("N31").Value = -10520.286
$ws.Range("M31").ClearContents()

$ws.Range("H34").Value = 9930.286
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 9930.286
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 9930.286
$ws.Range("N34").Value = -10334.286
$ws.Range("M34").ClearContents()

$ws.Range("H58").Value = 1688.3462
$ws.Range("I58").Value = 1246.4
$ws.Range("J58").Value = 2291
$ws.Range("K58").Value = 1246.4
$ws.Range("L58").Value = 2291
$ws.Range("M58").Value = -1043.4

$ws.Range("H122").Value = 992
$ws.Range("I122").Value = 992
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2976
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -526
$ws.Range("N122").ClearContents()

$ws.Range("H136").Value = 1688.3462
$ws.Range("I136").Value = 1246.4
$ws.Range("J136").Value = 2291
$ws.Range("K136").Value = 3739.2
$ws.Range("L136").Value = 6873
$ws.Range("M136").Value = -1189.2

$ws.Range("H138").Value = 59584.75
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 59584.75
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 59584.75
$ws.Range("N138").Value = -69864.75

$ws.Range("H139").Value = 29786.334
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 29786.334
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 29786.334
$ws.Range("N139").Value = -40066.334

$ws.Range("H140").Value = 23350.871
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 23350.871
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 23350.871
$ws.Range("N140").Value = -33710.871

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H31").Value = 800
$ws.Range("I31").Value = 800
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 2400
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -2112

$ws.Range("H56").Value = 4750
$ws.Range("I56").Value = 4750
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 4750
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = -4220

$ws.Range("H129").Value = 41669770
$ws.Range("I129").Value = 83336370
$ws.Range("J129").Value = 3166.5
$ws.Range("K129").Value = 250009110
$ws.Range("L129").Value = 9499.5
$ws.Range("M129").Value = -250004110
$ws.Range("N129").Value = -19499.5

$ws.Range("H133").Value = 49817.87
$ws.Range("I133").Value = 94164.63
$ws.Range("J133").Value = 9166.666999999999
$ws.Range("K133").Value = 282493.89
$ws.Range("L133").Value = 27500.001
$ws.Range("M133").Value = -277433.89
$ws.Range("N133").Value = -37620.001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4323620
$ws.Range("I122").Value = 4987638.5
$ws.Range("J122").Value = 7499.5
$ws.Range("K122").Value = 14962915.5
$ws.Range("L122").Value = 22498.5
$ws.Range("M122").Value = -14960465.5
$ws.Range("N122").Value = -27398.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 30304370
$ws.Range("I46").Value = 66667690
$ws.Range("J46").Value = 1603.6666
$ws.Range("K46").Value = 66667690
$ws.Range("L46").Value = 1603.6666
$ws.Range("M46").Value = -66667502
$ws.Range("N46").Value = -1979.6666

$ws.Range("H82").Value = 790962.3
$ws.Range("I82").Value = 1251073
$ws.Range("J82").Value = 177481.33
$ws.Range("K82").Value = 1251073
$ws.Range("L82").Value = 177481.33
$ws.Range("M82").Value = -1250712
$ws.Range("N82").Value = -178203.33

$ws.Range("H85").Value = 790962.3
$ws.Range("I85").Value = 1251073
$ws.Range("J85").Value = 177481.33
$ws.Range("K85").Value = 1251073
$ws.Range("L85").Value = 177481.33
$ws.Range("M85").Value = -1249825
$ws.Range("N85").Value = -179977.33

$ws.Range("H93").Value = 19240242
$ws.Range("I93").Value = 12706.277
$ws.Range("J93").Value = 62502200
$ws.Range("K93").Value = 12706.277
$ws.Range("L93").Value = 62502200
$ws.Range("M93").Value = -11458.277
$ws.Range("N93").Value = -62504696

$ws.Range("H122").Value = 8145407.5
$ws.Range("I122").Value = 11907845
$ws.Range("J122").Value = 2501751.2
$ws.Range("K122").Value = 35723535
$ws.Range("L122").Value = 7505253.600000001
$ws.Range("M122").Value = -35721085
$ws.Range("N122").Value = -7510153.600000001

$ws.Range("H136").Value = 7903.1943
$ws.Range("I136").Value = 6354.4614
$ws.Range("J136").Value = 11929.9
$ws.Range("K136").Value = 19063.3842
$ws.Range("L136").Value = 35789.7
$ws.Range("M136").Value = -16513.3842
$ws.Range("N136").Value = -40889.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1882.8572
$ws.Range("I81").Value = 1696.6666
$ws.Range("J81").Value = 3000
$ws.Range("K81").Value = 3393.3332
$ws.Range("L81").Value = 6000
$ws.Range("M81").Value = -2332.3332
$ws.Range("N81").Value = -8122

$ws.Range("H84").Value = 1882.8572
$ws.Range("I84").Value = 1696.6666
$ws.Range("J84").Value = 3000
$ws.Range("K84").Value = 16966.666
$ws.Range("L84").Value = 30000
$ws.Range("M84").Value = -11662.666
$ws.Range("N84").Value = -40608

$ws.Range("H96").Value = 2155.6155
$ws.Range("I96").Value = 1739.8
$ws.Range("J96").Value = 2415.5
$ws.Range("K96").Value = 1739.8
$ws.Range("L96").Value = 2415.5
$ws.Range("M96").Value = -366.8
$ws.Range("N96").Value = -5161.5

$ws.Range("H122").Value = 2387.5
$ws.Range("I122").Value = 1191.6666
$ws.Range("J122").Value = 5975
$ws.Range("K122").Value = 3574.9998
$ws.Range("L122").Value = 17925
$ws.Range("M122").Value = -1124.9998
$ws.Range("N122").Value = -22825

$ws.Range("H132").Value = 3064.9524
$ws.Range("I132").Value = 3386.1
$ws.Range("J132").Value = 2773
$ws.Range("K132").Value = 10158.3
$ws.Range("L132").Value = 8319
$ws.Range("M132").Value = -7628.299999999999
$ws.Range("N132").Value = -13379
